$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

# Correction to deal with nuclear in CES:
# Guarantee dispatch for additional non-dispatchable/renewable sources
# (onshore wind, solar PV, solar thermal, biomass, geothermal, offshore wind)
# by switching their base-year (B column) value from 0 to 1. All the other
# year columns (C:AK) are formulas referencing $B<row> and recalculate
# automatically.
$ws.Range("B6").Value = 1    # onshore wind
$ws.Range("B7").Value = 1    # solar PV
$ws.Range("B8").Value = 1    # solar thermal
$ws.Range("B9").Value = 1    # biomass
$ws.Range("B10").Value = 1   # geothermal
$ws.Range("B14").Value = 1   # offshore wind

# Make BGDPbES the active sheet/tab, with B11 selected, matching the
# author's last saved view state.
$ws.Activate()
$ws.Range("B11").Select()
